# Update statistical-errors table: new model columns (wB97M2, revDSD-PBEP86-D4)
# inserted, removed old DSD-PBEPBE/MGGA_MS2h columns, relabeled/reordered some
# "Mean ..." rows, inserted a new "Mean Thermochemistry" row, and refreshed all values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range('B1').Value = 'wB97M2'
$ws.Range('C1').Value = 'wB97M-V'
$ws.Range('D1').Value = 'CF22D'
$ws.Range('E1').Value = 'wB97X-V'
$ws.Range('F1').Value = 'revDSD-PBEP86-D4'
$ws.Range('G1').Value = 'M052X'
$ws.Range('H1').Value = 'M062X'
$ws.Range('I1').Value = 'M08HX'
$ws.Range('J1').Value = 'MN15'
$ws.Range('K1').Value = 'r2SCAN0'
$ws.Range('L1').Value = 'PW6B95'
$ws.Range('M1').Value = 'PBE0'
$ws.Range('N1').Value = 'CAMB3LYP'
$ws.Range('O1').Value = 'SOGGA11X'
$ws.Range('P1').Value = 'BMK'
$ws.Range('Q1').Value = 'B3LYP'

# Row 2
$ws.Range('A2').Value = 'Mean'
$ws.Range('B2').Value = 0.9485324850586064
$ws.Range('C2').Value = 1.075538502341131
$ws.Range('D2').Value = 1.24925710505858
$ws.Range('E2').Value = 1.312905121447046
$ws.Range('F2').Value = 1.521959467035516
$ws.Range('G2').Value = 1.941493346062861
$ws.Range('H2').Value = 1.971004445883499
$ws.Range('I2').Value = 1.993475061921084
$ws.Range('J2').Value = 2.009869356372007
$ws.Range('K2').Value = 2.811673997298874
$ws.Range('L2').Value = 3.277356932074
$ws.Range('M2').Value = 3.477682270471748
$ws.Range('N2').Value = 3.64800509733057
$ws.Range('O2').Value = 4.114856776341658
$ws.Range('P2').Value = 4.523620098569826
$ws.Range('Q2').Value = 5.691160285540097

# Row 3
$ws.Range('A3').Value = 'Mean Barrier Height'
$ws.Range('B3').Value = 0.6682301509999483
$ws.Range('C3').Value = 1.061577794867771
$ws.Range('D3').Value = 1.117046707592176
$ws.Range('E3').Value = 1.567957998053436
$ws.Range('F3').Value = 1.216549573109172
$ws.Range('G3').Value = 1.980175216602123
$ws.Range('H3').Value = 1.538446146087755
$ws.Range('I3').Value = 1.253832196919273
$ws.Range('J3').Value = 1.367874893291673
$ws.Range('K3').Value = 2.075873118633174
$ws.Range('L3').Value = 1.939300075785767
$ws.Range('M3').Value = 2.43049842142372
$ws.Range('N3').Value = 2.261613590144001
$ws.Range('O3').Value = 1.538501449436001
$ws.Range('P3').Value = 1.631565374781163
$ws.Range('Q3').Value = 3.292905240210722

# Row 4
$ws.Range('A4').Value = 'Mean Electric field'
$ws.Range('B4').Value = 1.03890582617713
$ws.Range('C4').Value = 1.807239366370509
$ws.Range('D4').Value = 1.981623492615641
$ws.Range('E4').Value = 1.040893097767846
$ws.Range('F4').Value = 1.887403100208598
$ws.Range('G4').Value = 1.062329342658534
$ws.Range('H4').Value = 1.020916802364601
$ws.Range('I4').Value = 1.205935272086035
$ws.Range('J4').Value = 1.53718374217268
$ws.Range('K4').Value = 1.116138939705384
$ws.Range('L4').Value = 1.472469255840157
$ws.Range('M4').Value = 1.568537354139071
$ws.Range('N4').Value = 1.298727770613418
$ws.Range('O4').Value = 1.089285966430315
$ws.Range('P4').Value = 1.045057997675376
$ws.Range('Q4').Value = 2.068887455683881

# Row 5
$ws.Range('A5').Value = 'Mean Frequency'
$ws.Range('B5').Value = 0.8761390189112278
$ws.Range('C5').Value = 0.9247261666421355
$ws.Range('D5').Value = 2.860531381523753
$ws.Range('E5').Value = 1.190705328840419
$ws.Range('F5').Value = 0.3264042308092379
$ws.Range('G5').Value = 1.322852217347581
$ws.Range('H5').Value = 1.013570648717671
$ws.Range('I5').Value = 1.25133300335017
$ws.Range('J5').Value = 1.762850698067449
$ws.Range('K5').Value = 1.666020535376949
$ws.Range('L5').Value = 1.455633822688355
$ws.Range('M5').Value = 0.7235033584737127
$ws.Range('N5').Value = 1.062130677838842
$ws.Range('O5').Value = 1.362824442117218
$ws.Range('P5').Value = 1.099713284993532
$ws.Range('Q5').Value = 1.060830648676765

# Row 6
$ws.Range('A6').Value = 'Mean Intramolecular Noncovalent'
$ws.Range('B6').Value = 1.021274463647786
$ws.Range('C6').Value = 1.070760594234186
$ws.Range('D6').Value = 1.254228493669825
$ws.Range('E6').Value = 0.9298746049867764
$ws.Range('F6').Value = 1.003151392433528
$ws.Range('G6').Value = 1.707520689303506
$ws.Range('H6').Value = 1.974323386417983
$ws.Range('I6').Value = 2.964708231009484
$ws.Range('J6').Value = 3.71778634826125
$ws.Range('K6').Value = 3.425001568048197
$ws.Range('L6').Value = 3.400001199884355
$ws.Range('M6').Value = 5.472002700319376
$ws.Range('N6').Value = 5.756046004437871
$ws.Range('O6').Value = 4.465077789117872
$ws.Range('P6').Value = 2.860608106593312
$ws.Range('Q6').Value = 8.06971402011504

# Row 7
$ws.Range('A7').Value = 'Mean Isomerization'
$ws.Range('B7').Value = 0.5554706594409674
$ws.Range('C7').Value = 1.30262562252924
$ws.Range('D7').Value = 1.22056851572857
$ws.Range('E7').Value = 1.903998258743017
$ws.Range('F7').Value = 3.066442516463451
$ws.Range('G7').Value = 1.761400181123495
$ws.Range('H7').Value = 1.608731078115534
$ws.Range('I7').Value = 1.594667840163088
$ws.Range('J7').Value = 1.76847405563724
$ws.Range('K7').Value = 2.762875133170416
$ws.Range('L7').Value = 3.012740638992309
$ws.Range('M7').Value = 3.435986542940126
$ws.Range('N7').Value = 3.603120007939581
$ws.Range('O7').Value = 3.252295269709014
$ws.Range('P7').Value = 2.340637216749883
$ws.Range('Q7').Value = 4.920660937042597

# Row 8
$ws.Range('A8').Value = 'Mean Noncovalent'
$ws.Range('B8').Value = 1.097451327629165
$ws.Range('C8').Value = 0.8907496469055344
$ws.Range('D8').Value = 1.408631313497661
$ws.Range('E8').Value = 0.9785908341633763
$ws.Range('F8').Value = 1.232059710847221
$ws.Range('G8').Value = 2.645796186753499
$ws.Range('H8').Value = 2.622894170617792
$ws.Range('I8').Value = 2.542658254902041
$ws.Range('J8').Value = 2.319516178483696
$ws.Range('K8').Value = 3.607798299111524
$ws.Range('L8').Value = 4.936263404742799
$ws.Range('M8').Value = 5.070317696035601
$ws.Range('N8').Value = 5.355664417772716
$ws.Range('O8').Value = 7.432988269066839
$ws.Range('P8').Value = 9.050913293739784
$ws.Range('Q8').Value = 9.16296128705798

# Row 9
$ws.Range('A9').Value = 'Mean Thermochemistry'
$ws.Range('B9').Value = 0.9246462968733978
$ws.Range('C9').Value = 1.109994036010521
$ws.Range('D9').Value = 1.05147331278402
$ws.Range('E9').Value = 1.567386731649206
$ws.Range('F9').Value = 1.653728971267002
$ws.Range('G9').Value = 1.463111649599972
$ws.Range('H9').Value = 1.578074208417043
$ws.Range('I9').Value = 1.571916769843446
$ws.Range('J9').Value = 1.78509715573366
$ws.Range('K9').Value = 2.525860194276377
$ws.Range('L9').Value = 2.621910765137486
$ws.Range('M9').Value = 2.353143933466784
$ws.Range('N9').Value = 2.5051831055744
$ws.Range('O9').Value = 2.391825476616795
$ws.Range('P9').Value = 2.597332979979269
$ws.Range('Q9').Value = 3.692115757406017

# Row 10
$ws.Range('A10').Value = 'Mean Transition Metal'
$ws.Range('B10').Value = 1.051872421684712
$ws.Range('C10').Value = 1.187662088716834
$ws.Range('D10').Value = 1.039248953079247
$ws.Range('E10').Value = 1.268111579774051
$ws.Range('F10').Value = 1.353032116155235
$ws.Range('G10').Value = 1.744404112116273
$ws.Range('H10').Value = 2.282591344538274
$ws.Range('I10').Value = 2.222966651198293
$ws.Range('J10').Value = 1.103916968429391
$ws.Range('K10').Value = 1.656858979878669
$ws.Range('L10').Value = 1.363581419509494
$ws.Range('M10').Value = 1.610199055411115
$ws.Range('N10').Value = 1.85042557165491
$ws.Range('O10').Value = 2.034809750464644
$ws.Range('P10').Value = 1.989480660640204
$ws.Range('Q10').Value = 2.128674061919511

# Row 10 is brand new; give the label cell the same bold/centered/bordered
# header style already used by the other column-A label cells (A2:A9).
$ws.Range('A2').Copy()
$ws.Range('A10').PasteSpecial(-4122)

[void]$ws.Range('A1').Select()
